$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("args")

# New command strings for column B, moving the previous argument strings to column D
$clientCmd = "java -jar Client.jar"
$serverCmd = "java-jar Server.jar"

# Row 28: ARGS Cliente
$old28 = $ws.Range("B28").Value2
$ws.Range("D28").Value = $old28
$ws.Range("B28").Value = $clientCmd

# Rows 29-37: ARGS nó x.y (server nodes)
for ($r = 29; $r -le 37; $r++) {
    $oldVal = $ws.Range("B$r").Value2
    $ws.Range("D$r").Value = $oldVal
    $ws.Range("B$r").Value = $serverCmd
}

# Copy formatting/style from column B to the newly used column D cells
$ws.Range("B28:B37").Copy() | Out-Null
$ws.Range("D28:D37").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Update the view state to match the authored workbook
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 25
$ws.Range("D39").Select() | Out-Null
